$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Diff: cell C10 on the sheet changes from 18 -> 1
$ws.Range("C10").Value = 1

